# Notas da prova 02
# Fill in the "Prova 02" column (C) for every student row, mirroring the
# per-student formula used to total the six grade components of Prova 01
# (column B), then drop a reviewer note on C16 and finish with the cell
# the author was last looking at (J5) selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows 5 and 13 are the "-" (no grade) rows, matching column B.
$ws.Range("C5").Value = "-"
$ws.Range("C13").Value = "-"

$ws.Range("C6").Formula  = "=15+15+20+10+17+0"
$ws.Range("C7").Formula  = "=15+15+20+4+20+10"
$ws.Range("C8").Formula  = "=15+12+17+10+20+18"
$ws.Range("C9").Formula  = "=15+15+12+10+20+20"
$ws.Range("C10").Formula = "=13+15+16+10+18+3"
$ws.Range("C11").Formula = "=15+15+20+5+20+8"
$ws.Range("C12").Formula = "=7+12+18+10+18+3"
$ws.Range("C14").Formula = "=14+15+12+8+20+5"
$ws.Range("C15").Formula = "=15+15+20+10+20+0"
$ws.Range("C16").Formula = "=15+15+15+10+20+2"
$ws.Range("C17").Formula = "=15+15+20+10+20+1"
$ws.Range("C18").Formula = "=15+15+20+8+20+12"
$ws.Range("C19").Formula = "=15+15+20+10+20+8"
$ws.Range("C20").Formula = "=15+15+20+8+20+3"
$ws.Range("C21").Formula = "=15+15+20+10+20+0"
$ws.Range("C22").Formula = "=15+15+20+10+20+5"
$ws.Range("C23").Formula = "=15+15+20+10+20+0"
$ws.Range("C24").Formula = "=2+10+20+10+5+0"
$ws.Range("C25").Formula = "=15+15+20+10+20+5"
$ws.Range("C26").Formula = "=15+15+20+3+18+2"
$ws.Range("C27").Formula = "=15+15+20+10+20+0"
$ws.Range("C28").Formula = "=15+15+20+8+20+12"
$ws.Range("C29").Formula = "=15+15+20+10+20+8"
$ws.Range("C30").Formula = "=15+15+20+10+20+0"
$ws.Range("C31").Formula = "=15+15+20+10+20+0"

# Reviewer comment on C16.
$ws.Range("C16").AddComment("Iago Augusto:`nO gato do token-ring tem razão")

# Leave the selection where the author ended up (also scrolls the view
# back to the top instead of the previous mid-sheet scroll position).
$ws.Range("J5").Select()
